$d = $word.ActiveDocument

$replacements = @(
    @("478÷7=", "672÷4="),
    @("186÷7=", "550÷2="),
    @("293÷5=", "681÷4="),
    @("753÷9=", "679÷3="),
    @("927÷2=", "270÷4="),
    @("896÷7=", "245÷9="),
    @("515÷8=", "470÷2="),
    @("976÷5=", "287÷3="),
    @("136÷6=", "558÷4="),
    @("676÷5=", "502÷2="),
    @("669÷4=", "891÷3="),
    @("212÷5=", "839÷7="),
    @("439÷6=", "588÷6="),
    @("510÷6=", "290÷8="),
    @("202÷9=", "198÷7="),
    @("432÷3=", "988÷6="),
    @("228÷8=", "560÷6="),
    @("563÷2=", "257÷9="),
    @("832÷8=", "532÷3="),
    @("879÷9=", "913÷8="),
    @("375÷3=", "503÷8="),
    @("437÷9=", "573÷2="),
    @("143÷6=", "419÷9="),
    @("444÷7=", "323÷4="),
    @("971÷8=", "594÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
